$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the D:E data range to Text format while writing values so that
# numeric-looking strings (e.g. "8.94", "0.0000297") are preserved exactly
# as text instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$origStyle = $dataRange.Style
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "96.921.10"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "3.695.72"
$ws.Range("E3").Value = "  +0.35%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "237.43"
$ws.Range("E5").Value = "  -2.64%  "

$ws.Range("D6").Value = "1.92"
$ws.Range("E6").Value = "  +2.09%  "

$ws.Range("D7").Value = "655.76"
$ws.Range("E7").Value = "  -1.82%  "

$ws.Range("D8").Value = "0.427"
$ws.Range("E8").Value = "  -0.36%  "

$ws.Range("E9").Value = "  -3.38%  "

$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D11").Value = "3.693.97"
$ws.Range("E11").Value = "  +0.40%  "

$ws.Range("D12").Value = "44.17"
$ws.Range("E12").Value = "  -2.86%  "

$ws.Range("E13").Value = "  +1.37%  "

$ws.Range("D14").Value = "0.0000297"
$ws.Range("E14").Value = "  +10.41%  "

$ws.Range("D15").Value = "6.77"
$ws.Range("E15").Value = "  +2.56%  "

$ws.Range("D16").Value = "4.385.37"
$ws.Range("E16").Value = "  +0.41%  "

$ws.Range("D17").Value = "96.724.57"
$ws.Range("E17").Value = "  +0.32%  "

$ws.Range("D18").Value = "8.94"
$ws.Range("E18").Value = "  -1.40%  "

$ws.Range("D19").Value = "3.709.30"
$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").Value = "13.01"
$ws.Range("E20").Value = "  +1.18%  "

$ws.Range("D21").Value = "18.63"
$ws.Range("E21").Value = "  +1.10%  "

$ws.Range("D22").Value = "0.509"
$ws.Range("E22").Value = "  -4.89%  "

$ws.Range("D23").Value = "522.12"
$ws.Range("E23").Value = "  +0.74%  "

$ws.Range("E24").Value = "  -1.95%  "

$ws.Range("E25").Value = "  +1.26%  "

$ws.Range("D26").Value = "6.91"
$ws.Range("E26").Value = "  -0.51%  "

$ws.Range("D27").Value = "101.74"
$ws.Range("E27").Value = "  -0.87%  "

$ws.Range("D28").Value = "0.190"
$ws.Range("E28").Value = "  +13.90%  "

$ws.Range("D29").Value = "13.41"
$ws.Range("E29").Value = "  +2.87%  "

$ws.Range("D30").Value = "12.28"
$ws.Range("E30").Value = "  +0.84%  "

$ws.Range("E31").Value = "  -1.46%  "

$ws.Range("E32").Value = "  +0.25%  "

$ws.Range("D33").Value = "0.188"
$ws.Range("E33").Value = "  +0.87%  "

$ws.Range("D34").Value = "1.86"
$ws.Range("E34").Value = "  +3.80%  "

$ws.Range("E35").Value = "  +0.23%  "

$ws.Range("D36").Value = "32.21"
$ws.Range("E36").Value = "  -2.06%  "

$ws.Range("D37").Value = "645.99"
$ws.Range("E37").Value = "  +4.45%  "

$ws.Range("D38").Value = "0.597"
$ws.Range("E38").Value = "  +1.95%  "

$ws.Range("D39").Value = "8.82"
$ws.Range("E39").Value = "  +0.90%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").Value = "6.85"
$ws.Range("E41").Value = "  +10.65%  "

$ws.Range("E42").Value = "  +5.19%  "

$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "0.160"
$ws.Range("E43").Value = "  +0.09%  "

$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").Value = "40.33"
$ws.Range("E44").Value = "  -5.03%  "

$ws.Range("D45").Value = "0.956"
$ws.Range("E45").Value = "  -0.20%  "

$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "0.452"
$ws.Range("E46").Value = "  +4.02%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0459"
$ws.Range("E47").Value = "  +0.13%  "

$ws.Range("D48").Value = "2.28"
$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("D49").Value = "23.63"
$ws.Range("E49").Value = "  +0.03%  "

$ws.Range("D50").Value = "8.54"
$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("D51").Value = "3.54"
$ws.Range("E51").Value = "  +1.21%  "

# Restore the original style/number format for the data range now that all
# text values have been written.
$dataRange.Style = $origStyle
